# New LEDs and mounting hardware
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: H1,H2 mounting hardware -> Wurth Elektronik terminal block
$ws.Range("D6").Value = "Wurth Elektronik"
$ws.Range("E6").Value = 7466303
$ws.Range("F6").Value = "Terminals WP-SMRA SMD Block 7mm"
$ws.Range("G6").Value = "SMD"
$ws.Range("H6").Value = "SMD"
$ws.Range("I6").Value = "Component should be mounted so screw hole faces left edge of board"

# Row 8: LEDs D8,D10,D5,D9,D11,D6,D7,D4 -> new Kingbright LED part
$ws.Range("E8").Value = "APPA3010SURCK"
$ws.Range("F8").Value = "Standard LEDs - SMD 3x1mm SMD RA RED "
$ws.Range("G8").Value = "SMD"

# Row 14: U4,U3 -> updated TI part / package
$ws.Range("E14").Value = "SN74HCS125QDRQ1"
$ws.Range("F14").Value = "Automotive Schmitt-trigger inputs quadruple bus buffer gates with 3-state outputs 14-SOIC -40 to 125 "

# View adjustments
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I19").Select()

$ws.Columns.Item(9).ColumnWidth = 57
